$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.578.86'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '2.289.28'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '95.67'
$ws.Range('E5').Value = '  +2.37%  '
$ws.Range('D6').Value = '267.51'
$ws.Range('E6').Value = '  -0.47%  '
$ws.Range('D7').Value = '0.623'
$ws.Range('E7').Value = '  -1.43%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('D9').Value = '0.606'
$ws.Range('E9').Value = '  -1.86%  '
$ws.Range('D10').Value = '45.83'
$ws.Range('E10').Value = '  +1.27%  '
$ws.Range('D11').Value = '0.0936'
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').Value = '7.75'
$ws.Range('E12').Value = '  -3.43%  '
$ws.Range('D13').Value = '0.106'
$ws.Range('E13').Value = '  +0.41%  '
$ws.Range('D14').Value = '2.630.71'
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('D15').Value = '15.06'
$ws.Range('E15').Value = '  -1.14%  '
$ws.Range('D16').Value = '0.848'
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('D17').Value = '2.284.15'
$ws.Range('E17').Value = '  -0.79%  '
$ws.Range('D18').Value = '43.560.57'
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('E19').Value = '  +2.58%  '
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').Value = '72.22'
$ws.Range('E21').Value = '  +1.75%  '
$ws.Range('D22').Value = '2.52'
$ws.Range('E22').Value = '  +10.65%  '
$ws.Range('D23').Value = "'231.90"
$ws.Range('E23').Value = '  -1.71%  '
$ws.Range('D24').Value = '9.15'
$ws.Range('E24').Value = '  -5.24%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').Value = '2.51'
$ws.Range('E26').Value = '  +0.45%  '
$ws.Range('D27').Value = '11.11'
$ws.Range('E27').Value = '  -0.70%  '
$ws.Range('E28').Value = '  +2.38%  '
$ws.Range('D29').Value = '40.55'
$ws.Range('E29').Value = '  +3.66%  '
$ws.Range('D30').Value = '2.22'
$ws.Range('E30').Value = '  -1.21%  '
$ws.Range('D31').Value = '175.58'
$ws.Range('E31').Value = '  +1.52%  '
$ws.Range('D32').Value = '21.83'
$ws.Range('E32').Value = '  -1.37%  '
$ws.Range('E33').Value = '  +1.45%  '
$ws.Range('E34').Value = '  -3.06%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  -1.21%  '
$ws.Range('E37').Value = '  +1.96%  '
$ws.Range('D38').Value = '4.34'
$ws.Range('E38').Value = '  -3.78%  '
$ws.Range('D39').Value = "'3.40"
$ws.Range('E39').Value = '  +0.73%  '
$ws.Range('E40').Value = '  +1.36%  '
$ws.Range('E41').Value = '  +0.89%  '
$ws.Range('D42').Value = '12.27'
$ws.Range('E42').Value = '  +0.34%  '
$ws.Range('D43').Value = "'65.50"
$ws.Range('E43').Value = '  +7.22%  '
$ws.Range('D44').Value = '1.35'
$ws.Range('E44').Value = '  +3.12%  '
$ws.Range('B45').Value = 'THORChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D45').Value = '5.22'
$ws.Range('E45').Value = '  -4.07%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '8.74'
$ws.Range('E46').Value = '  -1.20%  '
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('D48').Value = '97.24'
$ws.Range('E48').Value = '  -2.50%  '
$ws.Range('E49').Value = '  +0.55%  '
$ws.Range('D50').Value = '0.435'
$ws.Range('E50').Value = '  +1.61%  '
$ws.Range('D51').Value = '2.511.67'
$ws.Range('E51').Value = '  -0.08%  '
